# Regenerate save_data: update column G ("K") values for rows 2-30
# to reflect strikeouts (K) instead of the previous Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 6
    6  = 3
    7  = 7
    8  = 4
    9  = 7
    10 = 4
    11 = 3
    12 = 5
    13 = 3
    14 = 7
    15 = 6
    16 = 1
    17 = 2
    18 = 8
    19 = 5
    20 = 5
    21 = 3
    22 = 4
    23 = 4
    24 = 2
    25 = 3
    26 = 4
    27 = 3
    28 = 4
    29 = 1
    30 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
